$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("T2").Value = 144
$ws.Range("T3").Value = 555
$ws.Range("T5").Value = 720
$ws.Range("T6").Value = 203
